$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K6").Value = 3506.1426
$ws.Range("I6").Value = 1168.7142
$ws.Range("M6").Value = -3394.1426
$ws.Range("H6").Value = 1047.625
$ws.Range("N32").Value = -2283.125
$ws.Range("K32").Value = 1000
$ws.Range("L32").Value = 1631.125
$ws.Range("I32").Value = 1000
$ws.Range("J32").Value = 1631.125
$ws.Range("M32").Value = -674
$ws.Range("H32").Value = 1561
$ws.Range("N125").Value = -22332255
$ws.Range("L125").Value = 22327335
$ws.Range("J125").Value = 2480815
$ws.Range("H125").Value = 1860741.9
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K4").Value = 94.2
$ws.Range("I4").Value = 94.2
$ws.Range("M4").Value = 21.8
$ws.Range("H4").Value = 93.5
$ws.Range("K5").Value = 600
$ws.Range("I5").Value = 600
$ws.Range("M5").Value = -488
$ws.Range("H5").Value = 600
$ws.Range("N102").Value = -7044
$ws.Range("K102").Value = 1486.1538
$ws.Range("L102").Value = 3800
$ws.Range("I102").Value = 1486.1538
$ws.Range("J102").Value = 3800
$ws.Range("M102").Value = 135.8462
$ws.Range("H102").Value = 2128.889
$ws.Range("N122").Value = -10718.9998
$ws.Range("K122").Value = 4953
$ws.Range("L122").Value = 5818.9998
$ws.Range("I122").Value = 1651
$ws.Range("J122").Value = 1939.6666
$ws.Range("M122").Value = -2503
$ws.Range("H122").Value = 1696.579
$ws.Range("N135").Value = -61580.582
$ws.Range("L135").Value = 51440.582
$ws.Range("J135").Value = 51440.582
$ws.Range("H135").Value = 51440.582
$ws.Range("N139").Value = -56006.934
$ws.Range("K139").Value = 20000
$ws.Range("L139").Value = 45726.934
$ws.Range("I139").Value = 20000
$ws.Range("J139").Value = 45726.934
$ws.Range("M139").Value = -14860
$ws.Range("H139").Value = 42700.234
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K4").Value = 600
$ws.Range("I4").Value = 600
$ws.Range("M4").Value = -485
$ws.Range("H4").Value = 600
$ws.Range("N22").Value = -846
$ws.Range("K22").Value = 490.9565
$ws.Range("L22").Value = 500
$ws.Range("I22").Value = 490.9565
$ws.Range("J22").Value = 500
$ws.Range("M22").Value = -317.9565
$ws.Range("H22").Value = 491.33334
$ws.Range("N86").Value = -3845.8
$ws.Range("K86").Value = 1735.8422
$ws.Range("L86").Value = 1599.8
$ws.Range("I86").Value = 1735.8422
$ws.Range("J86").Value = 1599.8
$ws.Range("M86").Value = -612.8422
$ws.Range("H86").Value = 1729.04
$ws.Range("N89").Value = -19231
$ws.Range("K89").Value = 8679.210999999999
$ws.Range("L89").Value = 7999
$ws.Range("I89").Value = 1735.8422
$ws.Range("J89").Value = 1599.8
$ws.Range("M89").Value = -3063.210999999999
$ws.Range("H89").Value = 1729.04
$ws.Range("N94").Value = -2694
$ws.Range("K94").Value = 594.6
$ws.Range("L94").Value = 1792
$ws.Range("I94").Value = 594.6
$ws.Range("J94").Value = 1792
$ws.Range("M94").Value = -143.6
$ws.Range("H94").Value = 936.7143
$ws.Range("N134").Value = -120981
$ws.Range("K134").Value = 42073890
$ws.Range("L134").Value = 115911
$ws.Range("I134").Value = 14024630
$ws.Range("J134").Value = 38637
$ws.Range("M134").Value = -42071355
$ws.Range("H134").Value = 10185338
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N7").Value = -322
$ws.Range("K7").Value = 76.2
$ws.Range("L7").Value = 96
$ws.Range("I7").Value = 76.2
$ws.Range("J7").Value = 96
$ws.Range("M7").Value = 36.8
$ws.Range("H7").Value = 81.85714
$ws.Range("N22").Value = $null
$ws.Range("K22").Value = 199.1
$ws.Range("L22").Value = 0
$ws.Range("I22").Value = 199.1
$ws.Range("J22").Value = 0
$ws.Range("M22").Value = 150.9
$ws.Range("H22").Value = 199.1
$ws.Range("K31").Value = 3624459.8
$ws.Range("I31").Value = 3624459.8
$ws.Range("M31").Value = -3624164.8
$ws.Range("H31").Value = 2035816.5
$ws.Range("K34").Value = 3624459.8
$ws.Range("I34").Value = 3624459.8
$ws.Range("M34").Value = -3624257.8
$ws.Range("H34").Value = 2035816.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N134").Value = -24205.9995
$ws.Range("K134").Value = 2745
$ws.Range("L134").Value = 14065.9995
$ws.Range("I134").Value = 915
$ws.Range("J134").Value = 4688.6665
$ws.Range("M134").Value = 2325
$ws.Range("H134").Value = 1580.9412
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N70").Value = -28575168
$ws.Range("K70").Value = 2020392.8
$ws.Range("L70").Value = 28574628
$ws.Range("I70").Value = 2020392.8
$ws.Range("J70").Value = 28574628
$ws.Range("M70").Value = -2020122.8
$ws.Range("H70").Value = 5708481
$ws.Range("N73").Value = -28576500
$ws.Range("K73").Value = 2020392.8
$ws.Range("L73").Value = 28574628
$ws.Range("I73").Value = 2020392.8
$ws.Range("J73").Value = 28574628
$ws.Range("M73").Value = -2019456.8
$ws.Range("H73").Value = 5708481
$ws.Range("N80").Value = -16734.333
$ws.Range("K80").Value = 4051.4
$ws.Range("L80").Value = 14738.333
$ws.Range("I80").Value = 4051.4
$ws.Range("J80").Value = 14738.333
$ws.Range("M80").Value = -3053.4
$ws.Range("H80").Value = 7517.4326
$ws.Range("N83").Value = -83675.66500000001
$ws.Range("K83").Value = 20257
$ws.Range("L83").Value = 73691.66500000001
$ws.Range("I83").Value = 4051.4
$ws.Range("J83").Value = 14738.333
$ws.Range("M83").Value = -15265
$ws.Range("H83").Value = 7517.4326
$ws.Range("N122").Value = -55568710
$ws.Range("K122").Value = 157371.246
$ws.Range("L122").Value = 55563810
$ws.Range("I122").Value = 52457.082
$ws.Range("J122").Value = 18521270
$ws.Range("M122").Value = -154921.246
$ws.Range("H122").Value = 7967662.5
$ws.Range("N126").Value = -15095.9228
$ws.Range("K126").Value = 54301.00199999999
$ws.Range("L126").Value = 10155.9228
$ws.Range("I126").Value = 18100.334
$ws.Range("J126").Value = 3385.3076
$ws.Range("M126").Value = -51831.00199999999
$ws.Range("H126").Value = 8032.1577
$ws.Range("N132").Value = -20993201
$ws.Range("K132").Value = 19051123.5
$ws.Range("L132").Value = 20988141
$ws.Range("I132").Value = 6350374.5
$ws.Range("J132").Value = 6996047
$ws.Range("M132").Value = -19048593.5
$ws.Range("H132").Value = 6511792.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N7").Value = -3152.5715
$ws.Range("K7").Value = 1712.2222
$ws.Range("L7").Value = 2928.5715
$ws.Range("I7").Value = 1712.2222
$ws.Range("J7").Value = 2928.5715
$ws.Range("M7").Value = -1600.2222
$ws.Range("H7").Value = 2244.375
$ws.Range("N22").Value = -18521240
$ws.Range("K22").Value = 733.3333
$ws.Range("L22").Value = 18520650
$ws.Range("I22").Value = 733.3333
$ws.Range("J22").Value = 18520650
$ws.Range("M22").Value = -438.3333
$ws.Range("H22").Value = 15153392
$ws.Range("N27").Value = -18520864
$ws.Range("K27").Value = 733.3333
$ws.Range("L27").Value = 18520650
$ws.Range("I27").Value = 733.3333
$ws.Range("J27").Value = 18520650
$ws.Range("M27").Value = -626.3333
$ws.Range("H27").Value = 15153392
$ws.Range("N40").Value = -3522
$ws.Range("K40").Value = 2812.75
$ws.Range("L40").Value = 3250
$ws.Range("I40").Value = 2812.75
$ws.Range("J40").Value = 3250
$ws.Range("M40").Value = -2676.75
$ws.Range("H40").Value = 2900.2
$ws.Range("N46").Value = -55557406
$ws.Range("K46").Value = 746.2941
$ws.Range("L46").Value = 55557030
$ws.Range("I46").Value = 746.2941
$ws.Range("J46").Value = 55557030
$ws.Range("M46").Value = -558.2941
$ws.Range("H46").Value = 19231768
$ws.Range("N55").Value = -584.18182
$ws.Range("K55").Value = 62500076
$ws.Range("L55").Value = 238.18182
$ws.Range("I55").Value = 62500076
$ws.Range("J55").Value = 238.18182
$ws.Range("M55").Value = -62499903
$ws.Range("H55").Value = 16666861
$ws.Range("N122").Value = -300007900
$ws.Range("K122").Value = 1884298.2
$ws.Range("L122").Value = 300003000
$ws.Range("I122").Value = 628099.4
$ws.Range("J122").Value = 100001000
$ws.Range("M122").Value = -1881848.2
$ws.Range("H122").Value = 6148816
$ws.Range("N126").Value = -13725.7145
$ws.Range("K126").Value = 5136.6666
$ws.Range("L126").Value = 8785.7145
$ws.Range("I126").Value = 1712.2222
$ws.Range("J126").Value = 2928.5715
$ws.Range("M126").Value = -2666.6666
$ws.Range("H126").Value = 2244.375
$ws.Range("N132").Value = -17168.8181
$ws.Range("K132").Value = 35719008
$ws.Range("L132").Value = 12108.8181
$ws.Range("I132").Value = 11906336
$ws.Range("J132").Value = 4036.2727
$ws.Range("M132").Value = -35716478
$ws.Range("H132").Value = 6213931.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K122").Value = 3886.2
$ws.Range("I122").Value = 1295.4
$ws.Range("M122").Value = -1436.2
$ws.Range("H122").Value = 1914.6666
$ws.Range("N126").Value = -8465
$ws.Range("K126").Value = 75001974
$ws.Range("L126").Value = 3525
$ws.Range("I126").Value = 25000658
$ws.Range("J126").Value = 1175
$ws.Range("M126").Value = -74999504
$ws.Range("H126").Value = 17857948
$ws.Range("N132").Value = -5115331.699999999
$ws.Range("K132").Value = 11775.7779
$ws.Range("L132").Value = 5110271.699999999
$ws.Range("I132").Value = 3925.2593
$ws.Range("J132").Value = 1703423.9
$ws.Range("M132").Value = -9245.777900000001
$ws.Range("H132").Value = 869124.5600000001

Write-Output "Applied 246 cell updates across 8 sheets"